$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-range strings) ---
$ws.Range("A8").Value = "Volume 29   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/21/2022  Through  11/27/2022"

# --- Cells that require a style/number-format change (copy format from a same-style donor cell, then set value) ---
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G15").PasteSpecial(-4122)

$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("H15").PasteSpecial(-4122)

$ws.Range("I14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = 4

$ws.Range("K14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = 200

$ws.Range("I14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = 3

$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F30").PasteSpecial(-4122)

# --- Plain value updates (style unchanged) ---
$ws.Range("M15").Value = -28.571428571428
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -80
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -18.75
$ws.Range("I16").Value = 217
$ws.Range("J16").Value = 163
$ws.Range("K16").Value = 33.128834355828
$ws.Range("L16").Value = 2.358490566037
$ws.Range("M16").Value = 68.217054263565
$ws.Range("N16").Value = -82.556270096463
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 500
$ws.Range("F17").Value = 10
$ws.Range("H17").Value = -23.076923076923
$ws.Range("I17").Value = 155
$ws.Range("J17").Value = 145
$ws.Range("K17").Value = 6.896551724137
$ws.Range("L17").Value = 43.518518518518
$ws.Range("M17").Value = 74.157303370786
$ws.Range("N17").Value = -43.014705882352
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 75
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -20.833333333333
$ws.Range("I18").Value = 232
$ws.Range("J18").Value = 211
$ws.Range("K18").Value = 9.952606635071
$ws.Range("L18").Value = -36.956521739130
$ws.Range("M18").Value = 5.936073059360
$ws.Range("N18").Value = -91.479985310319
$ws.Range("C19").Value = 27
$ws.Range("D19").Value = 33
$ws.Range("E19").Value = -18.181818181818
$ws.Range("F19").Value = 122
$ws.Range("G19").Value = 137
$ws.Range("H19").Value = -10.948905109489
$ws.Range("I19").Value = 1579
$ws.Range("J19").Value = 1079
$ws.Range("K19").Value = 46.339202965709
$ws.Range("L19").Value = 28.792822185970
$ws.Range("M19").Value = 37.543554006968
$ws.Range("N19").Value = -53.667840375586
$ws.Range("C20").Value = 12
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 160
$ws.Range("I20").Value = 174
$ws.Range("J20").Value = 157
$ws.Range("K20").Value = 10.828025477707
$ws.Range("L20").Value = 55.357142857142
$ws.Range("M20").Value = 107.142857142857
$ws.Range("N20").Value = -94.476190476190
$ws.Range("C21").Value = 53
$ws.Range("D21").Value = 47
$ws.Range("E21").Value = 12.765957446808
$ws.Range("F21").Value = 190
$ws.Range("G21").Value = 200
$ws.Range("H21").Value = -5
$ws.Range("I21").Value = 2369
$ws.Range("J21").Value = 1769
$ws.Range("K21").Value = 33.917467495760
$ws.Range("L21").Value = 16.184404119666
$ws.Range("M21").Value = 40.676959619952
$ws.Range("N21").Value = -78.127596713138
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 44
$ws.Range("J22").Value = 17
$ws.Range("K22").Value = 158.823529411765
$ws.Range("L22").Value = 37.5
$ws.Range("M22").Value = 83.333333333333
$ws.Range("F23").Value = 5
$ws.Range("I23").Value = 36
$ws.Range("K23").Value = 20
$ws.Range("L23").Value = 16.129032258064
$ws.Range("M23").Value = 50
$ws.Range("C24").Value = 75
$ws.Range("D24").Value = 42
$ws.Range("E24").Value = 78.571428571428
$ws.Range("F24").Value = 356
$ws.Range("G24").Value = 186
$ws.Range("H24").Value = 91.397849462365
$ws.Range("I24").Value = 3684
$ws.Range("J24").Value = 2224
$ws.Range("K24").Value = 65.647482014388
$ws.Range("L24").Value = 53.308364544319
$ws.Range("M24").Value = 140.156453715776
$ws.Range("C25").Value = 4
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 25
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = 19.047619047619
$ws.Range("I25").Value = 344
$ws.Range("J25").Value = 309
$ws.Range("K25").Value = 11.326860841423
$ws.Range("L25").Value = 37.051792828685
$ws.Range("M25").Value = 7.165109034267
$ws.Range("J26").Value = 21
$ws.Range("K26").Value = -33.333333333333
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 75
$ws.Range("J27").Value = 86
$ws.Range("K27").Value = -12.790697674418
$ws.Range("L27").Value = 20.967741935483

$excel.CutCopyMode = 0
